$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1209.8
$ws.Cells.Item(6, 9).Value = 282.66666
$ws.Cells.Item(6, 10).Value = 2600.5
$ws.Cells.Item(6, 11).Value = 847.9999799999999
$ws.Cells.Item(6, 12).Value = 7801.5
$ws.Cells.Item(6, 13).Value = -735.9999799999999
$ws.Cells.Item(6, 14).Value = -8025.5

$ws.Cells.Item(21, 8).Value = 23286.334
$ws.Cells.Item(21, 9).Value = 23286.334
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 23286.334
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = -22818.334

$ws.Cells.Item(23, 8).Value = 23286.334
$ws.Cells.Item(23, 9).Value = 23286.334
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 23286.334
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = -23052.334

$ws.Cells.Item(137, 8).Value = 1630.4166
$ws.Cells.Item(137, 9).Value = 1161.0625
$ws.Cells.Item(137, 10).Value = 2569.125
$ws.Cells.Item(137, 11).Value = 3483.1875
$ws.Cells.Item(137, 12).Value = 7707.375
$ws.Cells.Item(137, 13).Value = -933.1875
$ws.Cells.Item(137, 14).Value = -12807.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 11628406
$ws.Cells.Item(2, 9).Value = 11628406
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 11628406
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -11628293

$ws.Cells.Item(32, 8).Value = 9049.846
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 9049.846
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 9049.846
$ws.Cells.Item(32, 14).Value = -9623.846
$ws.Cells.Item(32, 13).ClearContents()

$ws.Cells.Item(116, 8).Value = 11628406
$ws.Cells.Item(116, 9).Value = 11628406
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 11628406
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -11626112

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 11628406
$ws.Cells.Item(3, 9).Value = 11628406
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 11628406
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -11628292

$ws.Cells.Item(20, 8).Value = 2114.2727
$ws.Cells.Item(20, 9).Value = 2072.889
$ws.Cells.Item(20, 10).Value = 2300.5
$ws.Cells.Item(20, 11).Value = 2072.889
$ws.Cells.Item(20, 12).Value = 2300.5
$ws.Cells.Item(20, 13).Value = -1825.889
$ws.Cells.Item(20, 14).Value = -2794.5

$ws.Cells.Item(125, 8).Value = 30000
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 30000
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 30000
$ws.Cells.Item(125, 14).Value = -39840

$ws.Cells.Item(137, 8).Value = 58709.75
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 58709.75
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 58709.75
$ws.Cells.Item(137, 14).Value = -68909.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 613.3570999999999
$ws.Cells.Item(22, 9).Value = 286
$ws.Cells.Item(22, 10).Value = 1049.8334
$ws.Cells.Item(22, 11).Value = 286
$ws.Cells.Item(22, 12).Value = 1049.8334
$ws.Cells.Item(22, 13).Value = 64
$ws.Cells.Item(22, 14).Value = -1749.8334

$ws.Cells.Item(31, 8).Value = 3559
$ws.Cells.Item(31, 9).Value = 1856.4445
$ws.Cells.Item(31, 10).Value = 8666.666999999999
$ws.Cells.Item(31, 11).Value = 1856.4445
$ws.Cells.Item(31, 12).Value = 8666.666999999999
$ws.Cells.Item(31, 13).Value = -1561.4445
$ws.Cells.Item(31, 14).Value = -9256.666999999999

$ws.Cells.Item(34, 8).Value = 3559
$ws.Cells.Item(34, 9).Value = 1856.4445
$ws.Cells.Item(34, 10).Value = 8666.666999999999
$ws.Cells.Item(34, 11).Value = 1856.4445
$ws.Cells.Item(34, 12).Value = 8666.666999999999
$ws.Cells.Item(34, 13).Value = -1654.4445
$ws.Cells.Item(34, 14).Value = -9070.666999999999

$ws.Cells.Item(58, 8).Value = 2071736.4
$ws.Cells.Item(58, 9).Value = 3106625.5
$ws.Cells.Item(58, 10).Value = 1958.4286
$ws.Cells.Item(58, 11).Value = 3106625.5
$ws.Cells.Item(58, 12).Value = 1958.4286
$ws.Cells.Item(58, 13).Value = -3106422.5
$ws.Cells.Item(58, 14).Value = -2364.4286

$ws.Cells.Item(107, 8).Value = 380.15
$ws.Cells.Item(107, 9).Value = 333
$ws.Cells.Item(107, 10).Value = 647.3333
$ws.Cells.Item(107, 11).Value = 333
$ws.Cells.Item(107, 12).Value = 647.3333
$ws.Cells.Item(107, 13).Value = 1587
$ws.Cells.Item(107, 14).Value = -4487.3333

$ws.Cells.Item(132, 8).Value = 1649.5312
$ws.Cells.Item(132, 9).Value = 1198.2609
$ws.Cells.Item(132, 10).Value = 2802.7778
$ws.Cells.Item(132, 11).Value = 3594.7827
$ws.Cells.Item(132, 12).Value = 8408.3334
$ws.Cells.Item(132, 13).Value = -1064.7827
$ws.Cells.Item(132, 14).Value = -13468.3334

$ws.Cells.Item(134, 8).Value = 1604.7727
$ws.Cells.Item(134, 9).Value = 1324.1765
$ws.Cells.Item(134, 10).Value = 2558.8
$ws.Cells.Item(134, 11).Value = 3972.5295
$ws.Cells.Item(134, 12).Value = 7676.400000000001
$ws.Cells.Item(134, 13).Value = -1437.5295
$ws.Cells.Item(134, 14).Value = -12746.4

$ws.Cells.Item(136, 8).Value = 2071736.4
$ws.Cells.Item(136, 9).Value = 3106625.5
$ws.Cells.Item(136, 10).Value = 1958.4286
$ws.Cells.Item(136, 11).Value = 9319876.5
$ws.Cells.Item(136, 12).Value = 5875.2858
$ws.Cells.Item(136, 13).Value = -9317326.5
$ws.Cells.Item(136, 14).Value = -10975.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 713.4286
$ws.Cells.Item(5, 9).Value = 667.3333
$ws.Cells.Item(5, 10).Value = 990
$ws.Cells.Item(5, 11).Value = 2001.9999
$ws.Cells.Item(5, 12).Value = 2970
$ws.Cells.Item(5, 13).Value = -1889.9999
$ws.Cells.Item(5, 14).Value = -3194

$ws.Cells.Item(122, 8).Value = 722.0769
$ws.Cells.Item(122, 9).Value = 580
$ws.Cells.Item(122, 10).Value = 810.875
$ws.Cells.Item(122, 11).Value = 5220
$ws.Cells.Item(122, 12).Value = 7297.875
$ws.Cells.Item(122, 13).Value = -2770
$ws.Cells.Item(122, 14).Value = -12197.875

$ws.Cells.Item(131, 8).Value = 798.24
$ws.Cells.Item(131, 9).Value = 318.14285
$ws.Cells.Item(131, 10).Value = 834.37634
$ws.Cells.Item(131, 11).Value = 954.4285500000001
$ws.Cells.Item(131, 12).Value = 2503.12902
$ws.Cells.Item(131, 13).Value = 4085.57145
$ws.Cells.Item(131, 14).Value = -12583.12902

$ws.Cells.Item(135, 8).Value = 713.4286
$ws.Cells.Item(135, 9).Value = 667.3333
$ws.Cells.Item(135, 10).Value = 990
$ws.Cells.Item(135, 11).Value = 6005.9997
$ws.Cells.Item(135, 12).Value = 8910
$ws.Cells.Item(135, 13).Value = -3470.9997
$ws.Cells.Item(135, 14).Value = -13980

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 148.75
$ws.Cells.Item(2, 9).Value = 15
$ws.Cells.Item(2, 10).Value = 193.33333
$ws.Cells.Item(2, 11).Value = 15
$ws.Cells.Item(2, 12).Value = 193.33333
$ws.Cells.Item(2, 13).Value = 98
$ws.Cells.Item(2, 14).Value = -419.33333

$ws.Cells.Item(11, 8).Value = 6618872
$ws.Cells.Item(11, 9).Value = 8274927.5
$ws.Cells.Item(11, 10).Value = 546669.3
$ws.Cells.Item(11, 11).Value = 8274927.5
$ws.Cells.Item(11, 12).Value = 546669.3
$ws.Cells.Item(11, 13).Value = -8274788.5
$ws.Cells.Item(11, 14).Value = -546947.3

$ws.Cells.Item(57, 8).Value = 23266.666
$ws.Cells.Item(57, 9).Value = 9800
$ws.Cells.Item(57, 10).Value = 30000
$ws.Cells.Item(57, 11).Value = 9800
$ws.Cells.Item(57, 12).Value = 30000
$ws.Cells.Item(57, 13).Value = -8980
$ws.Cells.Item(57, 14).Value = -31640

$ws.Cells.Item(97, 8).Value = 2127.75
$ws.Cells.Item(97, 9).Value = 2500
$ws.Cells.Item(97, 10).Value = 1011
$ws.Cells.Item(97, 11).Value = 2500
$ws.Cells.Item(97, 12).Value = 1011
$ws.Cells.Item(97, 13).Value = -2004
$ws.Cells.Item(97, 14).Value = -2003

$ws.Cells.Item(102, 8).Value = 2171.8518
$ws.Cells.Item(102, 9).Value = 2078.3572
$ws.Cells.Item(102, 10).Value = 2272.5386
$ws.Cells.Item(102, 11).Value = 2078.3572
$ws.Cells.Item(102, 12).Value = 2272.5386
$ws.Cells.Item(102, 13).Value = -456.3571999999999
$ws.Cells.Item(102, 14).Value = -5516.5386

$ws.Cells.Item(113, 8).Value = 1549.8334
$ws.Cells.Item(113, 9).Value = 2000
$ws.Cells.Item(113, 10).Value = 1459.8
$ws.Cells.Item(113, 11).Value = 2000
$ws.Cells.Item(113, 12).Value = 1459.8
$ws.Cells.Item(113, 13).Value = 170
$ws.Cells.Item(113, 14).Value = -5799.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5236.9443
$ws.Cells.Item(7, 9).Value = 2226.9
$ws.Cells.Item(7, 10).Value = 8999.5
$ws.Cells.Item(7, 11).Value = 2226.9
$ws.Cells.Item(7, 12).Value = 8999.5
$ws.Cells.Item(7, 13).Value = -2114.9
$ws.Cells.Item(7, 14).Value = -9223.5

$ws.Cells.Item(16, 8).Value = 2737.0667
$ws.Cells.Item(16, 9).Value = 4016.3333
$ws.Cells.Item(16, 10).Value = 818.1667
$ws.Cells.Item(16, 11).Value = 4016.3333
$ws.Cells.Item(16, 12).Value = 818.1667
$ws.Cells.Item(16, 13).Value = -3846.3333
$ws.Cells.Item(16, 14).Value = -1158.1667

$ws.Cells.Item(68, 8).Value = 5791.6
$ws.Cells.Item(68, 9).Value = 6659.6665
$ws.Cells.Item(68, 10).Value = 4489.5
$ws.Cells.Item(68, 11).Value = 6659.6665
$ws.Cells.Item(68, 12).Value = 4489.5
$ws.Cells.Item(68, 13).Value = -5910.6665
$ws.Cells.Item(68, 14).Value = -5987.5

$ws.Cells.Item(71, 8).Value = 5791.6
$ws.Cells.Item(71, 9).Value = 6659.6665
$ws.Cells.Item(71, 10).Value = 4489.5
$ws.Cells.Item(71, 11).Value = 33298.3325
$ws.Cells.Item(71, 12).Value = 22447.5
$ws.Cells.Item(71, 13).Value = -29554.3325
$ws.Cells.Item(71, 14).Value = -29935.5

$ws.Cells.Item(122, 8).Value = 3859.2222
$ws.Cells.Item(122, 9).Value = 1804.6923
$ws.Cells.Item(122, 10).Value = 9201
$ws.Cells.Item(122, 11).Value = 5414.0769
$ws.Cells.Item(122, 12).Value = 27603
$ws.Cells.Item(122, 13).Value = -2964.0769
$ws.Cells.Item(122, 14).Value = -32503

$ws.Cells.Item(126, 8).Value = 5236.9443
$ws.Cells.Item(126, 9).Value = 2226.9
$ws.Cells.Item(126, 10).Value = 8999.5
$ws.Cells.Item(126, 11).Value = 6680.700000000001
$ws.Cells.Item(126, 12).Value = 26998.5
$ws.Cells.Item(126, 13).Value = -4210.700000000001
$ws.Cells.Item(126, 14).Value = -31938.5

$ws.Cells.Item(132, 8).Value = 1852.6875
$ws.Cells.Item(132, 9).Value = 1319.5555
$ws.Cells.Item(132, 10).Value = 2538.1428
$ws.Cells.Item(132, 11).Value = 3958.6665
$ws.Cells.Item(132, 12).Value = 7614.428400000001
$ws.Cells.Item(132, 13).Value = -1428.6665
$ws.Cells.Item(132, 14).Value = -12674.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3001.3333
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 3001.3333
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 3001.3333
$ws.Cells.Item(96, 14).Value = -5747.3333
